# Finished implementing the new tables for the public facing info pages.
#
# Inserts a new "Call of The Sea" specialty-shop item as row 27 on the
# "Items" sheet, pushing all the existing rows (27-72) down by one
# (28-73) while keeping their data/formatting intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 27 ("Pirate Leather
# Sleeves"). Excel shifts row 27..72 down to 28..73 and extends the
# used range to A1:BT73 automatically.
$ws.Rows.Item(27).Insert()

# Populate the newly inserted row 27 with the new item's data.
$ws.Range("A27").Value2 = 687
$ws.Range("C27").Value2 = "Call of The Sea"
$ws.Range("D27").Value2 = "spell-healing"
$ws.Range("F27").Value2 = "Pirate Lord Leather"
$ws.Range("G27").Value2 = "The call of the sea will heal the wounds of the desolate"

$ws.Range("J27").Value2 = 1000000
$ws.Range("L27").Value2 = 1000000000
$ws.Range("O27").Value2 = 0

$ws.Range("Q27").Value2 = 0.25
$ws.Range("X27").Value2 = 0.7
$ws.Range("Y27").Value2 = 0
$ws.Range("Z27").Value2 = 0.7

$ws.Range("AC27").Value2 = 1

$ws.Range("AI27").Value2 = 0
$ws.Range("AJ27").Value2 = 0
$ws.Range("AK27").Value2 = 0
$ws.Range("AL27").Value2 = 0
$ws.Range("AM27").Value2 = 0

$ws.Range("AS27").Value2 = 0
$ws.Range("AV27").Value2 = 0
$ws.Range("AX27").Value2 = 0
$ws.Range("AY27").Value2 = 0

$ws.Range("AZ27").Value2 = 1
$ws.Range("BA27").Value2 = 1
$ws.Range("BB27").Value2 = 0
$ws.Range("BC27").Value2 = 0
$ws.Range("BD27").Value2 = 0
$ws.Range("BE27").Value2 = 0
$ws.Range("BF27").Value2 = 0
$ws.Range("BG27").Value2 = 0

$ws.Range("BM27").Value2 = 0
$ws.Range("BN27").Value2 = 0
$ws.Range("BO27").Value2 = 0
$ws.Range("BP27").Value2 = 0
$ws.Range("BQ27").Value2 = 0

# Leave the selection where the author last left it.
$ws.Range("J62").Select()
